$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------
# Step 1: copy cell formatting (style only) into brand-new cells
# while source cells still carry their original (pre-edit) format.
# ---------------------------------------------------------------
$ws.Range("B18:C18").Copy()
$ws.Range("B17:C17").PasteSpecial(-4122)

$ws.Range("B21:C21").Copy()
$ws.Range("B22:C22").PasteSpecial(-4122)

$ws.Range("A22").Copy()
$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("A24").PasteSpecial(-4122)

$ws.Range("B23:C23").Copy()
$ws.Range("B25:C25").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------
# Step 2: overwrite cell text content for every row whose values
# actually change (row 10, rows 13-23) plus the two brand-new
# rows (24, 25).
# ---------------------------------------------------------------
# Row 10
$ws.Range('A10').Value = 'Objetivos:'
$ws.Range('B10').Value = 'Possibilitar ao estudante de Engenharia de Materiais o acesso a ferramentas computacionais modernas, de modo a que consiga descrever e quantificar conceitos vistos em outras disciplinas, como Ciência dos Materiais, Diagramas de Fases, Cinética de Transformação em Materiais, Termodinâmica, Propriedades Elétricas, Magnéticas, Térmicas e Ópticas, etc. Ao final do curso, o aluno será capaz de aplicar e entender resultados de simulações computacionais realistas aplicadas a diversas classes de materiais.'
$ws.Range('C10').Value = 'Possibilitar ao estudante de Engenharia de Materiais o acesso a ferramentas computacionais modernas, de modo a que consiga descrever e quantificar conceitos vistos em outras disciplinas, como Ciência dos Materiais, Diagramas de Fases, Cinética de Transformação em Materiais, Termodinâmica, Propriedades Elétricas, Magnéticas, Térmicas e Ópticas, etc. Ao final do curso, o aluno será capaz de aplicar e entender resultados de simulações computacionais realistas aplicadas a diversas classes de materiais.'

# Row 13
$ws.Range('B13').Value = '3480026 - João Paulo Pascon'
$ws.Range('C13').Value = '3480026 - João Paulo Pascon'

# Row 14
$ws.Range('B14').Value = '1176388 - Luiz Tadeu Fernandes Eleno'
$ws.Range('C14').Value = '1176388 - Luiz Tadeu Fernandes Eleno'

# Row 15
$ws.Range('A15').Value = 'Programa resumido:'
$ws.Range('B15').Value = 'Tratamento de imagens em materialografia; Ajuste de equações empíricas ; Potenciais interatômicos e dinâmica molecular clássica; Descrição da Cinética de nucleação e crescimento; Método dos Elementos Finitos; Métodos de Monte Carlo; Crescimento de grão; Cálculo de Diagramas de fases.'
$ws.Range('C15').Value = 'Tratamento de imagens em materialografia; Ajuste de equações empíricas ; Potenciais interatômicos e dinâmica molecular clássica; Descrição da Cinética de nucleação e crescimento; Método dos Elementos Finitos; Métodos de Monte Carlo; Crescimento de grão; Cálculo de Diagramas de fases.'

# Row 16
$ws.Range('A16').Value = 'Short syllabus:'
$ws.Range('B16').Value = 'Image processing in materialography; Adjusting empirical equations; Interatomic potentials and classical molecular dynamics; Description of nucleation and growth kinetics; Finite Element Method; Monte Carlo methods; Grain growth; Calculation of phase diagrams.'
$ws.Range('C16').Value = 'Image processing in materialography; Adjusting empirical equations; Interatomic potentials and classical molecular dynamics; Description of nucleation and growth kinetics; Finite Element Method; Monte Carlo methods; Grain growth; Calculation of phase diagrams.'

# Row 17
$ws.Range('A17').Value = 'Programa:'
$ws.Range('B17').Value = '- Tratamento de imagens: resolução, definição, contraste, saturação; uso de técnicas automatizadas de determinação de tamanho e distribuição de partículas.- Proposição e ajuste de equações empíricas a resultados de medidas experimentais: as diversas propostas de relações para a deformação plástica e encruamento.- Potenciais interatômicos e o método de dinâmica molecular clássica; simulação de solidificação de um metal puro.- Cinética de nucleação e crescimento: a equação de Johnson-Mehl-Avrami-Kolmogorov (JMAK) e sua aplicação computacional.- Elementos finitos: estudo do estado de tensão de materiais sob carregamentos mecânicos; simulação de transferência de calor em tratamentos térmicos.- Método de Monte Carlo aplicado à transição ferro-paramagnética e à cinética de crescimento de grão- Cálculo de diagramas de fases: curvas de energia livre, o método CALPHAD; Thermo-Calc e Dictra.'
$ws.Range('C17').Value = '- Tratamento de imagens: resolução, definição, contraste, saturação; uso de técnicas automatizadas de determinação de tamanho e distribuição de partículas.- Proposição e ajuste de equações empíricas a resultados de medidas experimentais: as diversas propostas de relações para a deformação plástica e encruamento.- Potenciais interatômicos e o método de dinâmica molecular clássica; simulação de solidificação de um metal puro.- Cinética de nucleação e crescimento: a equação de Johnson-Mehl-Avrami-Kolmogorov (JMAK) e sua aplicação computacional.- Elementos finitos: estudo do estado de tensão de materiais sob carregamentos mecânicos; simulação de transferência de calor em tratamentos térmicos.- Método de Monte Carlo aplicado à transição ferro-paramagnética e à cinética de crescimento de grão- Cálculo de diagramas de fases: curvas de energia livre, o método CALPHAD; Thermo-Calc e Dictra.'

# Row 18
$ws.Range('A18').Value = 'Syllabus:'
$ws.Range('B18').Value = '- Image treatment: resolution, definition, contrast, saturation; use of automated techniques for determining particle size and distribution.- Proposition and fit of empirical equations to results of experimental measures: the various proposals for relationships for plastic deformation and hardening.- Interatomic potentials and the classical molecular dynamics method; simulation of solidification of a pure metal.- Nucleation and growth kinetics: the Johnson-Mehl-Avrami-Kolmogorov (JMAK) equation and its computational application.- Finite element method: study of the stress state of materials under mechanical loads; simulation of heat transfer applied to heat treatments.- Monte Carlo method applied to the ferro-paramagnetic transition and to grain growth kinetics- Calculation of phase diagrams: free energy curves, the CALPHAD method; Thermo-Calc and Dictra.'
$ws.Range('C18').Value = '- Image treatment: resolution, definition, contrast, saturation; use of automated techniques for determining particle size and distribution.- Proposition and fit of empirical equations to results of experimental measures: the various proposals for relationships for plastic deformation and hardening.- Interatomic potentials and the classical molecular dynamics method; simulation of solidification of a pure metal.- Nucleation and growth kinetics: the Johnson-Mehl-Avrami-Kolmogorov (JMAK) equation and its computational application.- Finite element method: study of the stress state of materials under mechanical loads; simulation of heat transfer applied to heat treatments.- Monte Carlo method applied to the ferro-paramagnetic transition and to grain growth kinetics- Calculation of phase diagrams: free energy curves, the CALPHAD method; Thermo-Calc and Dictra.'

# Row 19
$ws.Range('A19').Value = 'Avaliação:'

# Row 20
$ws.Range('A20').Value = 'Método:'
$ws.Range('B20').Value = 'Aulas expositivas e em laboratório computacional, trabalhos e exercícios comentados. Trabalho baseado em Projeto'
$ws.Range('C20').Value = 'Aulas expositivas e em laboratório computacional, trabalhos e exercícios comentados. Trabalho baseado em Projeto'

# Row 21
$ws.Range('A21').Value = 'Critério:'
$ws.Range('B21').Value = 'Média aritmética de trabalhos propostos ao longo do curso (60%) e do Trabalho final em grupo (40%).'
$ws.Range('C21').Value = 'Média aritmética de trabalhos propostos ao longo do curso (60%) e do Trabalho final em grupo (40%).'

# Row 22
$ws.Range('A22').Value = 'Norma de recuperação:'
$ws.Range('B22').Value = 'Não haverá exame de recuperação.'
$ws.Range('C22').Value = 'Não haverá exame de recuperação.'

# Row 23
$ws.Range('A23').Value = 'Bibliografia:'
$ws.Range('B23').Value = '- Richard LESAR, Computational Materials Science – Fundamentals to Applications. MRS, 2013.- Rob Phillips, Crystals, Defects and Microstructures – Modelling across scales. Cambridge, 2001.- Artigos publicados em revistas como Computational Materials Science, Calphad, Journal of Alloys and Compounds, etc.'
$ws.Range('C23').Value = '- Richard LESAR, Computational Materials Science – Fundamentals to Applications. MRS, 2013.- Rob Phillips, Crystals, Defects and Microstructures – Modelling across scales. Cambridge, 2001.- Artigos publicados em revistas como Computational Materials Science, Calphad, Journal of Alloys and Compounds, etc.'

# Row 24
$ws.Range('A24').Value = 'Requisitos:'

# Row 25
$ws.Range('B25').Value = 'LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito fraco)
'
$ws.Range('C25').Value = 'LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito fraco)
'

# ---------------------------------------------------------------
# Step 3: clear cells that must become blank (old data from the
# pre-edit row alignment that doesn't exist in the new layout).
# ---------------------------------------------------------------
$ws.Range('A13').ClearContents()
$ws.Range('A14').ClearContents()
$ws.Range('B19').ClearContents()
$ws.Range('C19').ClearContents()

# ---------------------------------------------------------------
# Step 4: row heights
# ---------------------------------------------------------------
$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(14).AutoFit()
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(16).RowHeight = 60
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(18).RowHeight = 120
$ws.Rows.Item(19).AutoFit()
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 60
$ws.Rows.Item(23).RowHeight = 120
$ws.Rows.Item(24).AutoFit()
$ws.Rows.Item(25).RowHeight = 30
